$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column I (rows 2-5): values change from 4 to 5
$ws.Range("I2").Value = 5
$ws.Range("I3").Value = 5
$ws.Range("I4").Value = 5
$ws.Range("I5").Value = 5

# Add a new row 6 with training schedule data
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 6
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = -5
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 12
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = "train_dim2_1"

# Update selection to I7 (matches diff's selection change)
$ws.Range("I7").Select()
